# Apply the "Correct PDO sensitivity descriptions in Lecture 02 materials"
# edit: rename the single worksheet to "Example1" and add a second
# worksheet "Example2" containing a PDO (points-to-double-the-odds)
# sensitivity example with two scoring scenarios (PDO=20, PDO=40).

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Example1"

# Insert the new sheet right after Example1 so the tab order is
# Example1, Example2.
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Example2"

# --- Row 1: PDO header + scenario values -------------------------------
$ws2.Range("C1").Value = "PDO"
$ws2.Range("D1").Value = 20
$ws2.Range("E1").Value = 40

# --- Row 2: base score "A" ----------------------------------------------
$ws2.Range("C2").Value = "A"
$ws2.Range("D2").Value = 600
$ws2.Range("E2").Value = 600

# --- Row 3: slope "B" -----------------------------------------------------
$ws2.Range("C3").Value = "B"
$ws2.Range("D3").Formula = "=D1/LN(2)"
$ws2.Range("E3").Formula = "=E1/LN(2)"

# --- Row 5: column headers (write Score_PDO20/40 before "delta" so the
#     shared-string table is built in the same order as the target file)
$ws2.Range("A5").Value = "p"
$ws2.Range("B5").Value = "odds"
$ws2.Range("C5").Value = "ln(odds)"
$ws2.Range("D5").Value = "Score_PDO20"
$ws2.Range("E5").Value = "Score_PDO40"
$ws2.Range("G5").Value = "Score_PDO20"
$ws2.Range("H5").Value = "Score_PDO40"

# --- Row 4: "delta" headers over the G:H comparison columns --------------
$ws2.Range("G4").Value = "delta"
$ws2.Range("H4").Value = "delta"

# --- Rows 6-10: p / odds / ln(odds) / scores / deltas --------------------
$ws2.Range("A6").Value  = 0.01
$ws2.Range("A7").Value  = 0.02
$ws2.Range("A8").Value  = 0.03
$ws2.Range("A9").Value  = 0.04
$ws2.Range("A10").Value = 0.05

$ws2.Range("B6").Formula   = "=A6/(1-A6)"
$ws2.Range("B7:B10").Formula = "=A7/(1-A7)"

$ws2.Range("C6").Formula   = "=LN(B6)"
$ws2.Range("C7:C10").Formula = "=LN(B7)"

$ws2.Range("D6").Formula = "=D`$2-D`$3*`$C6"
$ws2.Range("E6").Formula = "=E`$2-E`$3*`$C6"
$ws2.Range("D7:E10").Formula = "=D`$2-D`$3*`$C7"

$ws2.Range("G7").Formula     = "=D7-D6"
$ws2.Range("H7").Formula     = "=E7-E6"
$ws2.Range("G8:G10").Formula = "=D8-D7"
$ws2.Range("H8:H10").Formula = "=E8-E7"

# Scores/deltas use the same thousands-separator number format as the
# Example1 sheet's Score column (reuses the existing style).
$ws2.Range("D6:E10").NumberFormat = "#,##0"
$ws2.Range("G7:H10").NumberFormat = "#,##0"

# Example2 becomes the active/selected sheet, with G8 selected.
$ws2.Activate()
$ws2.Range("G8").Select()
